# Applies the "Test Switch Logic; Continuous Logic with bugs" edit:
#  - Row 12 (Marine / Shield Defense / switch) gets a new companion row
#    inserted right after it for a new "Adrenaline" continuous skill.
#  - The old row 13 (Rockhead / Shocking) shifts down to row 14.
#  - A small vertical "detail form" for the new Adrenaline skill is added
#    at B19:C26.
#  - Column I (the "degree" column) switches to a Text number format on
#    the cells that hold data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row before the current row 13 (Rockhead/Shocking) ---
# so row12 (Marine/Shield Defense) is followed by the new Adrenaline row,
# and the old row13 becomes row14.
$ws.Rows.Item(13).Insert()

# --- 2. Fill in the new row 13: Marine / Adrenaline / continuous ---
# (column I / "degree" is filled in further down, after the detail form,
# to match the authoring order of the original edit)
$ws.Range("A13").Value = "Marine"
$ws.Range("B13").Value = "Adrenaline"
$ws.Range("C13").Value = "continuous"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "self"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = "uhp/speed/freq"

# copy row 12's formatting (fonts/alignment) down onto the new row 13 cells
$ws.Range("A12:H12").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Tweak existing cells around the new row ---
# E12 ("n" -> "self")
$ws.Range("E12").Value = "self"

# Column I switches to Text format on the header + existing data rows too.
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I14").NumberFormat = "@"

# --- 4. Add the new vertical "detail form" rows 19-26 describing the
#        freshly-added Adrenaline skill (row 13) ---
$ws.Range("B19").Value = "skill"
$ws.Range("C19").Value = "Adrenaline"

$ws.Range("B20").Value = "type"
$ws.Range("C20").Value = "continuous"

$ws.Range("B21").Value = "period"
$ws.Range("C21").Value = 2

$ws.Range("B22").Value = "tar_type"
$ws.Range("C22").Value = "self"

$ws.Range("B23").Value = "manacost"
$ws.Range("C23").Value = 0

$ws.Range("B24").Value = "cold_t"
$ws.Range("C24").Value = 4

$ws.Range("B25").Value = "attributes"
$ws.Range("C25").Value = "uhp/speed/freq"

$ws.Range("B26").Value = "degree"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "-20/1/2"

# Style the B column labels (bold/centered) like the other label column,
# and center the C column values.
$ws.Range("B11").Copy()
$ws.Range("B19:B26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C12").Copy()
$ws.Range("C19:C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column I (the "degree" column) uses a Text number format in the new
# layout. Row 13's "formula" is actually stored as literal text (not a
# real formula, hence the leading apostrophe to force text) - matching
# the bug referenced in the commit message.
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "'=F17"

# --- 5. Update the sheet view / selection to match the authored state ---
$ws.Range("I13").Select()
$ws.Application.ActiveWindow.ScrollRow = 4

# --- 6. Touch page setup (orientation) as in the source edit ---
$ws.PageSetup.Orientation = 1 # xlPortrait
